$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.336.37"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "2.045.52"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").Value = "'229.03"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'56.53"
$ws.Range("E8").Value = "  -3.51%  "

$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  -2.49%  "

$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("D12").Value = "2.350.41"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "'14.56"
$ws.Range("E13").Value = "  -3.67%  "

$ws.Range("D14").Value = "'20.55"
$ws.Range("E14").Value = "  -3.87%  "

$ws.Range("D15").Value = "'0.753"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").Value = "'5.27"
$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "2.043.60"
$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Value = "37.200.19"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D19").Value = "'5.97"
$ws.Range("E19").Value = "  -2.87%  "

$ws.Range("D20").Value = "'69.74"
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "'225.93"
$ws.Range("E22").Value = "  -1.96%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "'2.27"
$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("D26").Value = "'9.47"
$ws.Range("E26").Value = "  -3.95%  "

$ws.Range("D27").Value = "'167.78"
$ws.Range("E27").Value = "  -2.37%  "

$ws.Range("D28").Value = "'1.39"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("E29").Value = "  -6.59%  "

$ws.Range("D30").Value = "'18.88"
$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("D32").Value = "'4.51"
$ws.Range("E32").Value = "  -4.63%  "

$ws.Range("D33").Value = "'4.56"
$ws.Range("E33").Value = "  -2.34%  "

$ws.Range("D34").Value = "'0.0610"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'2.39"
$ws.Range("E35").Value = "  -4.27%  "

$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "'1.01"
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").Value = "'3.19"
$ws.Range("E38").Value = "  -6.64%  "

$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").Value = "'0.0220"
$ws.Range("E40").Value = "  -7.76%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "'16.94"
$ws.Range("E42").Value = "  +0.47%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.476.96"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("D44").Value = "'0.0939"
$ws.Range("E44").Value = "  -4.01%  "

$ws.Range("D45").Value = "'95.52"
$ws.Range("E45").Value = "  -7.04%  "

$ws.Range("D46").Value = "'1.15"

$ws.Range("E47").Value = "  -4.93%  "

$ws.Range("D48").Value = "'7.10"
$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("D49").Value = "'2.91"
$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").Value = "2.235.33"
$ws.Range("E50").Value = "  -1.98%  "

$ws.Range("D51").Value = "'3.62"
$ws.Range("E51").Value = "  -13.65%  "
